# Update the "N_ZnO" (K) column values per the plot cleanup / update commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2:K4").Value = 0.99299205054378525
$ws.Range("K5:K7").Value = 6.0389999881982703
$ws.Range("K8:K10").Value = 2.380480101313593
$ws.Range("K11:K13").Value = 4.0238998566627524
$ws.Range("K14:K16").Value = 2.02397991526365
$ws.Range("K17:K19").Value = 4.7801999125957479
$ws.Range("K20:K22").Value = 3.8102400576782189
$ws.Range("K23:K25").Value = 0.58100000293552856
$ws.Range("K26:K28").Value = 0.4130099979150289
$ws.Range("K29:K31").Value = 4.673900016403195
$ws.Range("K32:K34").Value = 9.7901999662399248
$ws.Range("K35:K37").Value = 1.9031999748229971
$ws.Range("K38:K40").Value = 8.3159999999999989
$ws.Range("K41:K43").Value = 3.2448000000000001
$ws.Range("K44:K46").Value = 4.3724999999999996
$ws.Range("K50:K52").Value = 4.3179999999999996
$ws.Range("K53:K55").Value = 1.8927
$ws.Range("K56:K58").Value = 1.0746
$ws.Range("K59:K61").Value = 1.3960999999999999
$ws.Range("K62:K64").Value = 0.95620000000000038
$ws.Range("K65:K67").Value = 8.0000000000000018
$ws.Range("K68:K70").Value = 1.7842
$ws.Range("K71:K73").Value = 3.7349999999999999
$ws.Range("K74:K76").Value = 1.2834000000000001
$ws.Range("K77:K79").Value = 4.862099999999999
$ws.Range("K80:K82").Value = 1.0289999999999999
$ws.Range("K83:K85").Value = 4.1216000000000008
$ws.Range("K95:K97").Value = 2.2334399999999999
$ws.Range("K101:K103").Value = 0.94300000000000017

# Restore the active-sheet selection to the full N_ZnO column, as left after the update.
$ws.Range("K1:K103").Select()
